$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": zero out this advisor's PORCELANATO (M) and
# NO RESURTIBLES (P) figures for client "BRITO CARDENAS RUTH CECILIA" (row 4),
# and update the corresponding "x de 8" progress counters in row 10.
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M4").Value = 0
$wsGrupo.Range("P4").Value = 0
$wsGrupo.Range("M10").Value = "1 de 8"
$wsGrupo.Range("P10").Value = "0 de 8"

# Sheet "VENTA MENSUAL": zero out septiembre (F) sales for the same client
# in row 4, and adjust the column total in row 10 accordingly.
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F4").Value = 0
$wsMensual.Range("F10").Value = 5372.02
